$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before the "/api/rooms" block (old row 9) to make
# room for the new "/api/signup" and "/api/signout" endpoint rows, pushing
# the rest of the table (rooms/tenants/transactions) down by two rows.
$ws.Range("A9:AI10").Insert()

# New row 8: /api/signup
$ws.Range("A8").Value = "/api/signup"
$ws.Range("B8").Value = "post/{username,password,email}"
$ws.Range("C8").Value = "{data:{token}}"

# New row 9: /api/signout
$ws.Range("A9").Value = "/api/signout"
$ws.Range("B9").Value = "delete"
$ws.Range("C9").Value = "{data:{message:”success}}"

# Rows 16/17 ("/api/tenants" and "/api/tenants/:id") had their return-type
# value sitting in column B; move it over to column C to match the other
# rows' layout.
$ws.Range("C16").Value = "{data:{tenants:[]}}"
$ws.Range("B16").Value = ""

$ws.Range("C17").Value = "{data:{tenant:{}}}"
$ws.Range("B17").Value = ""

# Append two new blank formatted rows at the bottom of the table (rows 37-38)
# matching the formatting of the existing trailing blank row.
$ws.Range("A36:AI36").Copy()
$ws.Range("A37:AI38").PasteSpecial(-4122)

# Leave the active selection on C11 (the "/api/rooms" return-type cell), as
# in the committed workbook.
$ws.Range("C11").Select()
